$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Vega Monumental Concepción"
$ws.Range("C14").Value = "Bíobío"
$ws.Range("D14").Value = 44799
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 100112022
$ws.Range("G14").Value = "Arveja Verde"
$ws.Range("H14").Value = "Perfection"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 40
$ws.Range("K14").Value = 20000
$ws.Range("L14").Value = 22000
$ws.Range("M14").Value = 21000
$ws.Range("N14").Value = "`$/malla 25 kilos"
$ws.Range("O14").Value = "Provincia de Huasco"
$ws.Range("P14").Value = 840
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
